$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = 45854.64411551836
$question = "q01_sessions_avg_per_week"

$donorIds = @(
    "0ce5dd49",
    "2c1001cb",
    "37cc37bf",
    "43faa0b9",
    "4abe3e88",
    "50164f59",
    "5cf70f79",
    "5da96769",
    "6ca3e2f6",
    "790a4fcb",
    "802cc63a",
    "85c3ea4d",
    "942dfafb",
    "9bc6ba8c",
    "a2d65af2",
    "a46f1771",
    "ad58f9da",
    "c7d9a301",
    "ce8732ff",
    "d6f1d567",
    "da9326c9",
    "e09ca7bf",
    "ef53a641"
)
$avgSessions = @(
    7.475247524752476,
    7.405797101449275,
    2.881578947368421,
    4.642857142857143,
    2.836065573770492,
    5.769230769230769,
    2.571428571428572,
    1.5,
    2.5,
    1.666666666666667,
    2.910714285714286,
    9.839285714285714,
    8.027027027027026,
    1,
    2.5,
    3.28125,
    2.39622641509434,
    3.848484848484849,
    2.985714285714286,
    1.727272727272727,
    2.857142857142857,
    8.027027027027026,
    7.475247524752476
)
$categories = @(
    "6-10",
    "6-10",
    "1-2",
    "3-5",
    "1-2",
    "3-5",
    "1-2",
    "1-2",
    "1-2",
    "1-2",
    "1-2",
    "6-10",
    "6-10",
    "1-2",
    "1-2",
    "3-5",
    "1-2",
    "3-5",
    "1-2",
    "1-2",
    "1-2",
    "6-10",
    "6-10"
)

for ($i = 0; $i -lt $donorIds.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $donorIds[$i]
    $ws.Cells.Item($row, 2).Value = $avgSessions[$i]
    $ws.Cells.Item($row, 3).Value = $categories[$i]
    $ws.Cells.Item($row, 4).Value = $question
    $ws.Cells.Item($row, 5).Value = $timestamp
    $ws.Cells.Item($row, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
